$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on Price (D) and Volume (E) columns for the rows
# we touch, so numeric-looking strings (e.g. "1.001") are stored as text
# exactly like the source inlineStr cells, not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '22.997.53'
$ws.Range("E2").Value = '  -3.91%  '

$ws.Range("D3").Value = '1.600.81'
$ws.Range("E3").Value = '  -3.02%  '

$ws.Range("E4").Value = '  +0.20%  '

$ws.Range("D5").Value = '1.001'
$ws.Range("E5").Value = '  +0.13%  '

$ws.Range("D6").Value = '300.73'
$ws.Range("E6").Value = '  -3.29%  '

$ws.Range("D7").Value = '0.3776'
$ws.Range("E7").Value = '  -3.06%  '

$ws.Range("D8").Value = '0.3623'
$ws.Range("E8").Value = '  -5.71%  '

$ws.Range("D9").Value = '49.86'
$ws.Range("E9").Value = '  -2.38%  '

$ws.Range("D10").Value = '1.258'
$ws.Range("E10").Value = '  -6.55%  '

$ws.Range("E11").Value = '  +0.19%  '

$ws.Range("D12").Value = '0.08115'
$ws.Range("E12").Value = '  -3.94%  '

$ws.Range("D13").Value = '22.79'
$ws.Range("E13").Value = '  -4.83%  '

$ws.Range("D14").Value = '6.591'
$ws.Range("E14").Value = '  -6.25%  '

$ws.Range("D15").Value = '7.359'
$ws.Range("E15").Value = '  -6.88%  '

$ws.Range("D16").Value = '0.00001239'
$ws.Range("E16").Value = '  -6.06%  '

$ws.Range("D17").Value = '1.598.65'
$ws.Range("E17").Value = '  -3.17%  '

$ws.Range("D18").Value = '92.10'
$ws.Range("E18").Value = '  -2.02%  '

$ws.Range("D19").Value = '0.06869'
$ws.Range("E19").Value = '  -1.42%  '

$ws.Range("D20").Value = '18.17'
$ws.Range("E20").Value = '  -7.19%  '

$ws.Range("D21").Value = '6.565'
$ws.Range("E21").Value = '  -5.39%  '

$ws.Range("D22").Value = '0.5560'
$ws.Range("E22").Value = '  -6.29%  '

$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.09%  '

$ws.Range("E24").Value = '  -3.61%  '

$ws.Range("D25").Value = '23.008.08'
$ws.Range("E25").Value = '  -3.87%  '

$ws.Range("D26").Value = '2.364'
$ws.Range("E26").Value = '  -2.96%  '

$ws.Range("D27").Value = '2.802'
$ws.Range("E27").Value = '  -4.01%  '

$ws.Range("D28").Value = '21.06'
$ws.Range("E28").Value = '  -4.18%  '

$ws.Range("D29").Value = '150.93'
$ws.Range("E29").Value = '  -2.16%  '

$ws.Range("D30").Value = '5.267'
$ws.Range("E30").Value = '  -2.19%  '

$ws.Range("D31").Value = '133.33'
$ws.Range("E31").Value = '  -2.87%  '

$ws.Range("E32").Value = '  -6.84%  '

$ws.Range("D33").Value = '6.786'
$ws.Range("E33").Value = '  -12.23%  '

$ws.Range("D34").Value = '1.778.44'
$ws.Range("E34").Value = '  -2.97%  '

$ws.Range("D35").Value = '0.9604'
$ws.Range("E35").Value = '  -3.21%  '

$ws.Range("D36").Value = '0.07625'

$ws.Range("D37").Value = '10.33'
$ws.Range("E37").Value = '  -1.58%  '

$ws.Range("D38").Value = '6.254'
$ws.Range("E38").Value = '  -6.90%  '

$ws.Range("D39").Value = '0.02699'
$ws.Range("E39").Value = '  -8.02%  '

$ws.Range("D40").Value = '0.2531'
$ws.Range("E40").Value = '  -5.98%  '

$ws.Range("D41").Value = '0.08867'
$ws.Range("E41").Value = '  -2.74%  '

$ws.Range("D42").Value = '1.366'
$ws.Range("E42").Value = '  -3.92%  '

$ws.Range("D43").Value = '0.7052'
$ws.Range("E43").Value = '  -6.77%  '

$ws.Range("D44").Value = '12.50'
$ws.Range("E44").Value = '  -6.96%  '

$ws.Range("D45").Value = '15.24'
$ws.Range("E45").Value = '  -8.48%  '

$ws.Range("D46").Value = '0.6616'
$ws.Range("E46").Value = '  -4.63%  '

$ws.Range("B47").Value = 'Frax'
$ws.Range("C47").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  +0.04%  '

$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.308'
$ws.Range("E48").Value = '  -5.67%  '

$ws.Range("D49").Value = '3.983'
$ws.Range("E49").Value = '  -2.85%  '

$ws.Range("D50").Value = '131.90'
$ws.Range("E50").Value = '  -1.73%  '

$ws.Range("D51").Value = '0.07903'
$ws.Range("E51").Value = '  -4.51%  '

# Restore the default "Normal" style so no spurious per-cell style/number
# format is left behind (matches original cells which carry no explicit
# style index).
$ws.Range("D2:E51").Style = "Normal"
